$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text tweaks on existing labels (column A) ---
$ws.Range("A5").Value = "Processing of measured variable"
$ws.Range("A7").Value = "Heath state- Health indicator mapping, fk"
$ws.Range("A8").Value = "Health indicator - Measured variable mapping, hk"

# --- New "He 2012" data column (column E) ---
$ws.Range("E2").Value  = "Spiral bevel gear degradation (pitting)"
$ws.Range("E3").Value  = "Oil debris"
$ws.Range("E4").Value  = "Oil debris, Acceleration"
$ws.Range("E5").Value  = "One dimensional transition function using whitening transform"
$ws.Range("E7").Value  = "Direct"
$ws.Range("E8").Value  = " Data driven Double exponential smoothing model"
$ws.Range("E10").Value = "N/A This was buildt on data (ARIMA)"
$ws.Range("E11").Value = "Particle Filter with l-step ahead estimator"

# --- Reviewer note, flagged with the built-in "Bad" cell style ---
$ws.Range("E14").Value = "Check if this is summarized in lit review"
$ws.Range("E14").Style = "Bad"

# --- Update view selection to match the author's last position ---
$null = $ws.Range("F18").Select()
